$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: move the "Shina" translation from C1 (D1 is a duplicate) into B1,
# then clear out the old C1/D1 cells entirely.
$shina = $ws.Range("C1").Value2
$ws.Range("B1").Value = $shina
$ws.Range("C1:D1").ClearContents()

# Rows 2-36: mirror column A's shared string into column B (new translation
# column) for every data row.
for ($r = 2; $r -le 36; $r++) {
    $srcCell = $ws.Cells.Item($r, 1)
    $dstCell = $ws.Cells.Item($r, 2)
    $dstCell.Value = $srcCell.Value2
}
